{"js": "// Append a new development-diary row to the end of the (only) table in the\n// document, matching the row immediately above it for formatting (borders /\n// shading) and filling in the new session's data.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst date = \"28/01/2022\";\nconst timeSpent = \"3 Hours\";\nconst segment = \"Simulation \\u2013 Objective 3\";\nconst notes =\n  \"Added small UI improvements as well as a new console and debug tab for \" +\n  \"the province viewer. The console can be used to activate debug mode \" +\n  \"(Currently by typing \\u201cDEBUG\\u201d) allowing a user access to a tab \" +\n  \"that displays internal information such as IDs of a province and its \" +\n  \"culture. This console will later be expanded to allow the use of the \" +\n  \"soon-to-be-implemented actions performable by a nation.\";\n\ntable.addRows(\"End\", 1, [[date, timeSpent, segment, notes]]);\nawait context.sync();\n", "ps1": "# Append a new development-diary row to the end of the (only) table in the\n# document, matching the row immediately above it for formatting (borders /\n# shading) and filling in the new session's data.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newRow = $t.Rows.Add()\n$rowIndex = $newRow.Index\n\n$date = \"28/01/2022\"\n$timeSpent = \"3 Hours\"\n$segment = \"Simulation \" + [char]0x2013 + \" Objective 3\"\n$notes = \"Added small UI improvements as well as a new console and debug tab for the province viewer. The console can be used to activate debug mode (Currently by typing \" + [char]0x201C + \"DEBUG\" + [char]0x201D + \") allowing a user access to a tab that displays internal information such as IDs of a province and its culture. This console will later be expanded to allow the use of the soon-to-be-implemented actions performable by a nation.\"\n\n$t.Cell($rowIndex, 1).Range.Text = $date\n$t.Cell($rowIndex, 2).Range.Text = $timeSpent\n$t.Cell($rowIndex, 3).Range.Text = $segment\n$t.Cell($rowIndex, 4).Range.Text = $notes\n"}
